$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "compent type" column (I) -------------------------------------
# Header cell for the new column.
$ws.Range("I1").Value = "compent type"

# Widen column I to fit the new, longer header text.
$ws.Columns.Item(9).ColumnWidth = 36.8340625

# --- Column H ("deactivate" = "yes") cleanup ---------------------------
# These rows no longer carry a "yes" marker in column H.
$rowsToClear = @(3,4,6,7,8,9,10,11,12,13)
foreach ($r in $rowsToClear) {
    $ws.Range("H$r").Clear()
}

# Row 14 now carries the "yes" marker that used to live elsewhere.
$ws.Range("H14").Value = "yes"

# --- Fill in missing id numbers for the trailing rows ------------------
$ws.Range("A38").Value = 36
$ws.Range("A39").Value = 37
$ws.Range("A40").Value = 38

# --- Restore the cursor / selection -------------------------------------
$ws.Range("I14").Select() | Out-Null
